$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.155.54'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '3.090.61'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.45'
$ws.Range("E5").Value = '  +2.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '619.71'
$ws.Range("E6").Value = '  -1.83%  '
$ws.Range("E7").Value = '  +9.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.368'
$ws.Range("E8").Value = '  +2.02%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").Value = '3.089.53'
$ws.Range("E10").Value = '  -0.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.752'
$ws.Range("E11").Value = '  +5.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.201'
$ws.Range("E12").Value = '  +3.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000250'
$ws.Range("E13").Value = '  +2.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.55'
$ws.Range("E14").Value = '  -3.18%  '
$ws.Range("D15").Value = '91.188.31'
$ws.Range("E15").Value = '  +0.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.48'
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("D17").Value = '3.669.21'
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").Value = '3.102.31'
$ws.Range("E18").Value = '  -1.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.67'
$ws.Range("E19").Value = '  -3.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.60'
$ws.Range("E20").Value = '  +3.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000212'
$ws.Range("E21").Value = '  +2.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.83'
$ws.Range("E22").Value = '  +4.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '445.05'
$ws.Range("E23").Value = '  +0.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.15'
$ws.Range("E24").Value = '  +2.16%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '91.22'
$ws.Range("E25").Value = '  +3.28%  '
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.63'
$ws.Range("E26").Value = '  -1.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.91'
$ws.Range("E27").Value = '  -5.50%  '
$ws.Range("D28").Value = '3.259.10'
$ws.Range("E28").Value = '  -1.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E30").Value = '  +17.30%  '
$ws.Range("E31").Value = '  +26.77%  '
$ws.Range("E32").Value = '  -2.86%  '
$ws.Range("E33").Value = '  +15.15%  '
$ws.Range("E34").Value = '  +0.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.111'
$ws.Range("E35").Value = '  +32.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.75'
$ws.Range("E36").Value = '  +8.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.54'
$ws.Range("E37").Value = '  +0.80%  '
$ws.Range("E38").Value = '  +28.65%  '
$ws.Range("B39").Value = 'PancakeSwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.92'
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '495.16'
$ws.Range("E40").Value = '  -2.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.62'
$ws.Range("E41").Value = '  -4.38%  '
$ws.Range("E42").Value = '  +1.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.421'
$ws.Range("E43").Value = '  +2.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.14'
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.90'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.694'
$ws.Range("E47").Value = '  +1.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '153.88'
$ws.Range("E48").Value = '  +1.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.48'
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.34'
$ws.Range("E50").Value = '  -0.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.19'
$ws.Range("E51").Value = '  -1.81%  '
